$d = $word.ActiveDocument
$d.Content.Find.Execute("6x6, 1.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6x6, 0.05", 2)
